# Apply the cryptocurrency price/volume updates described by the commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.941.47"
$ws.Range("E2").Value = "  +2.33%  "
$ws.Range("D3").Value = "1.580.58"
$ws.Range("E3").Value = "  +1.86%  "
$ws.Range("E4").Value = "  -0.48%  "
$cell = $ws.Range("D5")
$origStyle = $cell.Style
$cell.Value = "'211.72"
$cell.Style = $origStyle
$ws.Range("E5").Value = "  +1.22%  "
$ws.Range("E6").Value = "  +7.48%  "
$ws.Range("E7").Value = "  -0.58%  "
$cell = $ws.Range("D8")
$origStyle = $cell.Style
$cell.Value = "'25.39"
$cell.Style = $origStyle
$ws.Range("E8").Value = "  +8.39%  "
$cell = $ws.Range("D9")
$origStyle = $cell.Style
$cell.Value = "'0.248"
$cell.Style = $origStyle
$ws.Range("E9").Value = "  +2.86%  "
$ws.Range("E10").Value = "  +1.37%  "
$ws.Range("E11").Value = "  +1.31%  "
$ws.Range("D12").Value = "1.803.99"
$ws.Range("E12").Value = "  +1.78%  "
$ws.Range("D13").Value = "1.555.15"
$ws.Range("E13").Value = "  +0.26%  "
$ws.Range("D14").Value = "28.896.15"
$ws.Range("E14").Value = "  +2.18%  "
$ws.Range("E15").Value = "  +2.39%  "
$ws.Range("E16").Value = "  +1.63%  "
$cell = $ws.Range("D17")
$origStyle = $cell.Style
$cell.Value = "'62.22"
$cell.Style = $origStyle
$ws.Range("E17").Value = "  +2.94%  "
$cell = $ws.Range("D18")
$origStyle = $cell.Style
$cell.Value = "'233.42"
$cell.Style = $origStyle
$ws.Range("E18").Value = "  +3.13%  "
$cell = $ws.Range("D19")
$origStyle = $cell.Style
$cell.Value = "'7.43"
$cell.Style = $origStyle
$ws.Range("E19").Value = "  +1.82%  "
$ws.Range("D20").Value = "0.0₃0692"
$ws.Range("E20").Value = "  +2.73%  "
$cell = $ws.Range("D21")
$origStyle = $cell.Style
$cell.Value = "'0.997"
$cell.Style = $origStyle
$ws.Range("E21").Value = "  -0.42%  "
$ws.Range("E22").Value = "  +1.81%  "
$cell = $ws.Range("D23")
$origStyle = $cell.Style
$cell.Value = "'9.17"
$cell.Style = $origStyle
$ws.Range("E23").Value = "  +4.10%  "
$cell = $ws.Range("D24")
$origStyle = $cell.Style
$cell.Value = "'2.10"
$cell.Style = $origStyle
$ws.Range("E24").Value = "  +4.57%  "
$cell = $ws.Range("D25")
$origStyle = $cell.Style
$cell.Value = "'152.46"
$cell.Style = $origStyle
$ws.Range("E25").Value = "  +3.09%  "
$ws.Range("E26").Value = "  +4.59%  "
$cell = $ws.Range("D27")
$origStyle = $cell.Style
$cell.Value = "'15.01"
$cell.Style = $origStyle
$ws.Range("E27").Value = "  +1.50%  "
$cell = $ws.Range("D28")
$origStyle = $cell.Style
$cell.Value = "'6.33"
$cell.Style = $origStyle
$ws.Range("E28").Value = "  +1.96%  "
$ws.Range("E29").Value = "  -0.49%  "
$cell = $ws.Range("D30")
$origStyle = $cell.Style
$cell.Value = "'0.0464"
$cell.Style = $origStyle
$ws.Range("E30").Value = "  -0.56%  "
$ws.Range("E31").Value = "  +0.34%  "
$cell = $ws.Range("D32")
$origStyle = $cell.Style
$cell.Value = "'3.21"
$cell.Style = $origStyle
$ws.Range("E32").Value = "  +1.36%  "
$ws.Range("D33").Value = "1.420.64"
$ws.Range("E33").Value = "  +2.59%  "
$ws.Range("E34").Value = "  -0.82%  "
$cell = $ws.Range("D35")
$origStyle = $cell.Style
$cell.Value = "'1.04"
$cell.Style = $origStyle
$ws.Range("E35").Value = "  -1.23%  "
$cell = $ws.Range("D36")
$origStyle = $cell.Style
$cell.Value = "'1.51"
$cell.Style = $origStyle
$ws.Range("E36").Value = "  +0.63%  "
$cell = $ws.Range("D37")
$origStyle = $cell.Style
$cell.Value = "'2.74"
$cell.Style = $origStyle
$ws.Range("E37").Value = "  +6.71%  "
$cell = $ws.Range("D38")
$origStyle = $cell.Style
$cell.Value = "'2.29"
$cell.Style = $origStyle
$ws.Range("E38").Value = "  -1.99%  "
$ws.Range("E39").Value = "  +1.01%  "
$cell = $ws.Range("D40")
$origStyle = $cell.Style
$cell.Value = "'0.525"
$cell.Style = $origStyle
$ws.Range("E40").Value = "  +2.74%  "
$ws.Range("E41").Value = "  +0.88%  "
$ws.Range("E42").Value = "  -0.48%  "
$cell = $ws.Range("D43")
$origStyle = $cell.Style
$cell.Value = "'0.785"
$cell.Style = $origStyle
$ws.Range("E43").Value = "  +1.25%  "
$cell = $ws.Range("D44")
$origStyle = $cell.Style
$cell.Value = "'0.0460"
$cell.Style = $origStyle
$ws.Range("E44").Value = "  -0.96%  "
$cell = $ws.Range("D45")
$origStyle = $cell.Style
$cell.Value = "'64.61"
$cell.Style = $origStyle
$ws.Range("E45").Value = "  +4.57%  "
$cell = $ws.Range("D46")
$origStyle = $cell.Style
$cell.Value = "'5.31"
$cell.Style = $origStyle
$ws.Range("E46").Value = "  -1.96%  "
$ws.Range("D47").Value = "1.716.45"
$ws.Range("E47").Value = "  +1.86%  "
$ws.Range("B48").Value = "BitcoinSV"
$ws.Range("C48").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$cell = $ws.Range("D48")
$origStyle = $cell.Style
$cell.Value = "'43.79"
$cell.Style = $origStyle
$ws.Range("E48").Value = "  +4.98%  "
$ws.Range("B49").Value = "WEMIXToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$cell = $ws.Range("D49")
$origStyle = $cell.Style
$cell.Value = "'0.839"
$cell.Style = $origStyle
$ws.Range("E49").Value = "  -7.48%  "
$ws.Range("B50").Value = "Quant"
$ws.Range("C50").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$cell = $ws.Range("D50")
$origStyle = $cell.Style
$cell.Value = "'85.37"
$cell.Style = $origStyle
$ws.Range("E50").Value = "  -0.02%  "
$cell = $ws.Range("D51")
$origStyle = $cell.Style
$cell.Value = "'0.0512"
$cell.Style = $origStyle
$ws.Range("E51").Value = "  +0.70%  "
